# Add a new "Users" worksheet at the end of the workbook, populate it with
# the header row + one sample data row, size the columns, and mark the
# numeric-looking text values as intentionally-text (numberStoredAsText).

$wb = $excel.ActiveWorkbook

# Append the new sheet after the last existing sheet so it lands at the end
# (Worksheets.Add() with no args inserts at the front, which is not what we want).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Users"

# Column widths (character units) matching the template's other sheets.
$colWidths = @(15.83203125, 15.83203125, 15.83203125, 30.83203125, 15.83203125, 12.83203125, 25.83203125, 30.83203125, 15.83203125, 10.83203125, 10.83203125, 15.83203125)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i] - (5/6)
}

$headers = @("userId","firstName","lastName","email","phone","role","organization","address","city","postcode","isActive","registrationDate")

# Columns that hold numeric/date-looking text which must stay text
# (phone number with leading zero, and a dd/mm/yyyy date string).
$textColumns = @(5, 12)
foreach ($col in $textColumns) {
    $ws.Cells.Item(2, $col).NumberFormat = "@"
}

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @("USER001","Alex","Smith","alex.smith@example.com","07123456789","user","Example Org","10 Example Street","Newcastle","NE1 1AA","Yes","01/01/2025")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# Flag the whole used range as "number stored as text" so Excel doesn't
# show the green-triangle warning, matching the other sheets in this workbook.
$ws.Range("A1:L2").ErrorCheckingOptions.NumberAsText = $true
